# TRIFFHIR-39: Import value sets from excel files
#
# Collapses the two-sheet "Valuesets" + "Concepts" workbook into a single
# "Concepts" sheet that also carries the value-set ID/Name, by deleting the
# old "Valuesets" sheet and inserting two new leading columns (ID, Name) into
# "Concepts", populated from the value-set URL already present in each row.

$wb = $excel.ActiveWorkbook

# Drop the separate "Valuesets" listing sheet - its Name/ID columns move
# into "Concepts" itself.
$wsValuesets = $wb.Worksheets.Item("Valuesets")
$wsValuesets.Delete() | Out-Null

$ws = $wb.Worksheets.Item("Concepts")

# Make room for the new ID/Name columns at the front; this shifts the
# existing ValueSet URL / Code / Display / Code System URL columns from
# A:D to C:F, carrying over their formatting (header bold, URL column
# fill style, widths) untouched.
$ws.Range("A:B").Insert() | Out-Null

# Header row
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "URL"
$ws.Cells.Item(1, 6).Value = "System"

# Row 2 & 3 both reference ValueSet v1 -> "vs1-test" / "Value Set 1"
$ws.Cells.Item(2, 1).Value = "vs1-test"
$ws.Cells.Item(2, 2).Value = "Value Set 1"
$ws.Cells.Item(3, 1).Value = "vs1-test"
$ws.Cells.Item(3, 2).Value = "Value Set 1"

# Row 4 references ValueSet v2 -> "vs2-test" / "Value Set 2"
$ws.Cells.Item(4, 1).Value = "vs2-test"
$ws.Cells.Item(4, 2).Value = "Value Set 2"

# Match the new ID column's width to its content.
$ws.Columns.Item(2).ColumnWidth = 10.0221354167

$ws.Range("F2").Select() | Out-Null

Write-Output "done"
